# Journal de travail - ajout des entrees de la Semaine 2 (fin) / Semaine 3
# (creation de la page d'administration)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Les lignes 28 a 36 etaient vides (mise en forme "hors tableau", styles 12-17).
# On copie la mise en forme de la derniere ligne remplie (27, styles 6/7/8/9/10/5)
# vers les lignes a completer, puis on saisit les nouvelles donnees.
$ws.Range("A27:F27").Copy()
$ws.Range("A28:F36").PasteSpecial(-4122)

# Ligne 28
$ws.Range("A28").Value = 45058
$ws.Range("B28").Value = 2
$ws.Range("C28").Value = 0.5
$ws.Range("D28").Value = "Implémentation"
$ws.Range("E28").Value = "Avancer sur la page des commandes du client"

# Ligne 29
$ws.Range("A29").Value = 45058
$ws.Range("B29").Value = 2
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = "Implémentation"
$ws.Range("E29").Value = "Héberger une partie du site sur SwissCenter"

# Ligne 30
$ws.Range("A30").Value = 45058
$ws.Range("B30").Value = 2
$ws.Range("C30").Value = 0.75
$ws.Range("D30").Value = "Analyse"
$ws.Range("E30").Value = "Srpint Review"

# Ligne 31
$ws.Range("A31").Value = 45058
$ws.Range("B31").Value = 2
$ws.Range("C31").Value = 0.75
$ws.Range("D31").Value = "Documentation"
$ws.Range("E31").Value = "Ajouter des stratgies de test dans le rappot"

# Ligne 32
$ws.Range("A32").Value = 45061
$ws.Range("B32").Value = 3
$ws.Range("C32").Value = 2.25
$ws.Range("D32").Value = "Implémentation"
$ws.Range("E32").Value = "Améliorer l'aspect graphique du site"

# Ligne 33 (texte plus long => hauteur de ligne agrandie)
$ws.Range("A33").Value = 45061
$ws.Range("B33").Value = 3
$ws.Range("C33").Value = 1.5
$ws.Range("D33").Value = "Implémentation"
$ws.Range("E33").Value = "Envoyer les données du panier vers la page des commandes"
$ws.Rows("33").RowHeight = 30

# Ligne 34
$ws.Range("A34").Value = 45061
$ws.Range("B34").Value = 3
$ws.Range("C34").Value = 1.5
$ws.Range("D34").Value = "Implémentation"
$ws.Range("E34").Value = "Créer la page administrateur et afficher les biscuits"

# Ligne 35
$ws.Range("A35").Value = 45061
$ws.Range("B35").Value = 3
$ws.Range("C35").Value = 1.5
$ws.Range("D35").Value = "Implémentation"
$ws.Range("E35").Value = "Coder l'ajout de produit et la modification"

# Ligne 36
$ws.Range("A36").Value = 45062
$ws.Range("B36").Value = 3
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = "Implémentation"
$ws.Range("E36").Value = "Finaliser la modification des biscuits"

# Mise a jour de la vue de la feuille : defilement + cellule active
$ws.Range("A37").Select()
